$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.824.02"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "2.305.09"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.95%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "2.656.97"
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").Value = "2.302.38"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.822"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "46.750.05"
$ws.Range("E18").Value = "  +6.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +22.52%  "
$ws.Range("D20").Value = "0.0₃0945"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.53%  "
$ws.Range("E24").Value = "  +5.95%  "
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "44.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.93%  "
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0802"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.18%  "
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E38").Value = "  +9.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0305"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.58%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.862.28"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.97%  "
$ws.Range("E47").Value = "  +9.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.70%  "
